# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates the "K" column (column G) values for each data row on Sheet1
# to reflect the recalculated strike-count (K) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of worksheet row number -> new value for column G ("K")
$newK = @{
    2  = 3
    3  = 1
    4  = 2
    5  = 2
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 0
    14 = 1
    15 = 1
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 0
    21 = 0
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 1
    27 = 2
    29 = 1
    30 = 1
    31 = 1
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 1
    39 = 2
    40 = 1
    41 = 1
    43 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
